$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.116.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.466.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.507'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.477.18'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0967'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  +4.42%  '
$ws.Range("E13").Value = '  +1.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.900.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.114.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.470.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("E19").Value = '  +6.21%  '
$ws.Range("E20").Value = '  +4.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '317.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("E23").Value = '  +7.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.410'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.161'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.574.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.83%  '
$ws.Range("E29").Value = '  +7.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0781'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '147.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.98%  '
$ws.Range("E34").Value = '  +4.05%  '
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("E36").Value = '  +8.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.71'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.857'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.78%  '
$ws.Range("E39").Value = '  +2.82%  '
$ws.Range("E40").Value = '  +7.01%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0550'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.602'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("E44").Value = '  +6.37%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.28%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '260.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +10.73%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0921'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.72%  '
$ws.Range("E49").Value = '  +3.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.870.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.25%  '
